$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.139.18"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "1.567.69"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'208.10"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'22.11"
$ws.Range("E8").Value = "  +4.20%  "
$ws.Range("D9").Value = "'0.250"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").Value = "'0.0588"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").Value = "'0.0859"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").Value = "1.792.03"
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("D13").Value = "1.569.01"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").Value = "27.107.95"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "'62.07"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "'219.19"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("D19").Value = ("{0}{1}{2}" -f "0.0", [char]0x2083, "0698")
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").Value = "'7.35"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("D23").Value = "'9.30"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "'154.19"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").Value = "'6.63"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("D28").Value = "'1.01"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("D30").Value = "'0.0471"
$ws.Range("E30").Value = "  +3.06%  "
$ws.Range("D31").Value = "'1.11"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "1.449.54"
$ws.Range("E33").Value = "  +6.33%  "
$ws.Range("D34").Value = "'3.08"
$ws.Range("E34").Value = "  +5.07%  "
$ws.Range("E35").Value = "  +4.08%  "
$ws.Range("D36").Value = "'0.967"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("D38").Value = "'0.0166"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("D39").Value = "'0.526"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").Value = "'0.815"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("D41").Value = "'5.75"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("E43").Value = "  +4.52%  "
$ws.Range("D44").Value = "'0.989"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'64.67"
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("E46").Value = "  +2.44%  "
$ws.Range("D47").Value = "1.705.47"
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("D48").Value = "'86.90"
$ws.Range("E48").Value = "  +3.34%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = ("{0}{1}{2}" -f "0.0", [char]0x2086, "0103")
$ws.Range("E49").Value = "  +5.76%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0525"
$ws.Range("E50").Value = "  +3.11%  "
$ws.Range("D51").Value = "'0.0966"
$ws.Range("E51").Value = "  +2.68%  "
